$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Row 32 ---
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A32:E32").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(32,1).Value2 = ' SCRIPT/G01P05A/um2101.ssb'
$ws.Cells.Item(32,2).Value2 = 338
$ws.Cells.Item(32,3).Value2 = ' Go talk to [CS:N]Torkoal[CR],\nthe town elder.'
$ws.Cells.Item(32,4).Value2 = ' Поговорите с [CS:N]Торкоалом[CR],\nгородским старейшиной.'
$ws.Cells.Item(32,5).Value2 = ' Ðïãïâïñéóå ò [CS:N]Óïñëïàìïí[CR],\nãïñïäòëéí òóàñåêšéîïê.'
$ws.Rows.Item(32).RowHeight = 57.6
$excel.CutCopyMode = $false

# --- Row 33 ---
$ws.Range("B2:E2").Copy() | Out-Null
$ws.Range("B33:E33").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(33,2).Value2 = 341
$ws.Cells.Item(33,3).Value2 = ' If anyone knows about the\n[CS:P]Hidden Land[CR], it\''d be him! ♪'
$ws.Cells.Item(33,4).Value2 = ' Если кто и знает про [CS:P]Сокрытые\nЗемли[CR], то это он! ♪'
$ws.Cells.Item(33,5).Value2 = ' Åòìé ëóï é èîàåó ðñï [CS:P]Òïëñúóúå\nÈåíìé[CR], óï üóï ïî! ♪'
$ws.Rows.Item(33).RowHeight = 31.8
$excel.CutCopyMode = $false

# --- Row 34 ---
$ws.Range("B2:E2").Copy() | Out-Null
$ws.Range("B34:E34").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(34,2).Value2 = 344
$ws.Cells.Item(34,3).Value2 = ' You should put his years of\nwisdom to use! ♪'
$ws.Cells.Item(34,4).Value2 = ' Вам пригодится его многолетняя\nмудрость! ♪'
$ws.Cells.Item(34,5).Value2 = ' Âàí ðñéãïäéóòÿ åãï íîïãïìåóîÿÿ\níôäñïòóû! ♪'
$ws.Rows.Item(34).RowHeight = 21.6
$excel.CutCopyMode = $false

# --- Row 35 ---
$ws.Range("B2:E2").Copy() | Out-Null
$ws.Range("B35:E35").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(35,2).Value2 = 313
$ws.Cells.Item(35,3).Value2 = ' You should ask [CS:N]Torkoal[CR], the\ntown elder.'
$ws.Cells.Item(35,4).Value2 = ' Вам нужно расспросить [CS:N]Торкоала[CR],\nгородского старейшину.'
$ws.Cells.Item(35,5).Value2 = ' Âàí îôçîï ñàòòðñïòéóû [CS:N]Óïñëïàìà[CR],\nãïñïäòëïãï òóàñåêšéîô.'
$ws.Rows.Item(35).RowHeight = 31.8
$excel.CutCopyMode = $false

# --- Row 36 ---
$ws.Range("B2:E2").Copy() | Out-Null
$ws.Range("B36:E36").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(36,2).Value2 = 316
$ws.Cells.Item(36,3).Value2 = ' He may know something about\nthe [CS:P]Hidden Land[CR]! ♪'
$ws.Cells.Item(36,4).Value2 = ' Он может знать что-нибудь о\n[CS:P]Сокрытых Землях[CR]! ♪'
$ws.Cells.Item(36,5).Value2 = ' Ïî íïçåó èîàóû œóï-îéáôäû ï\n[CS:P]Òïëñúóúö Èåíìÿö[CR]! ♪'
$ws.Rows.Item(36).RowHeight = 21.6
$excel.CutCopyMode = $false

# --- Row 37 ---
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A37:E37").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(37,2).Value2 = 319
$ws.Cells.Item(37,3).Value2 = ' You should put his years of\nwisdom to use! ♪'
$ws.Cells.Item(37,4).Value2 = ' Вам пригодится его многолетняя\nмудрость! ♪'
$ws.Cells.Item(37,5).Value2 = ' Âàí ðñéãïäéóòÿ åãï íîïãïìåóîÿÿ\níôäñïòóû! ♪'
$ws.Rows.Item(37).RowHeight = 21.6
$excel.CutCopyMode = $false

# --- Row 38 ---
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A38:E38").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(38,1).Value2 = 'SCRIPT/G01P05A/um2201.ssb'
$ws.Cells.Item(38,2).Value2 = 291
$ws.Cells.Item(38,3).Value2 = ' We\''re going to [CS:P]Brine Cave[CR]\ntomorrow! ♪'
$ws.Cells.Item(38,4).Value2 = ' Завтра мы отправимся в [CS:P]Пещеру\nу Моря[CR]! ♪'
$ws.Cells.Item(38,5).Value2 = ' Èàâóñà íú ïóðñàâéíòÿ â [CS:P]Ðåþåñô\nô Íïñÿ[CR]! ♪'
$ws.Rows.Item(38).RowHeight = 43.2
$excel.CutCopyMode = $false

# --- Row 39 ---
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A39:E39").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(39,1).Value2 = 'SCRIPT/G01P05A/um2202.ssb'
$ws.Cells.Item(39,2).Value2 = 294
$ws.Cells.Item(39,3).Value2 = ' Go get ready while you have\nthe chance! ♪'
$ws.Cells.Item(39,4).Value2 = ' Подготовьтесь, пока у вас есть\nтакая возможность!♪'
$ws.Cells.Item(39,5).Value2 = ' Ðïäãïóïâûóåòû, ðïëà ô âàò åòóû\nóàëàÿ âïèíïçîïòóû! ♪'
$ws.Rows.Item(39).RowHeight = 43.2
$excel.CutCopyMode = $false

# --- Row 40 ---
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A40:E40").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(40,1).Value2 = ' SCRIPT/G01P05A/um2401.ssb'
$ws.Cells.Item(40,2).Value2 = 269
$ws.Cells.Item(40,3).Value2 = ' Ooooh, Team [team:]! You\ncan do it! ♪'
$ws.Cells.Item(40,4).Value2 = ' Ооооо, Команда [team:]!\nУ вас всё получится! ♪'
$ws.Cells.Item(40,5).Value2 = ' Ïïïïï, Ëïíàîäà [team:]!\nÔ âàò âòæ ðïìôœéóòÿ! ♪'
$ws.Rows.Item(40).RowHeight = 28.8
$excel.CutCopyMode = $false

# --- Row 41 ---
$ws.Range("B2:E2").Copy() | Out-Null
$ws.Range("B41:E41").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(41,2).Value2 = 272
$ws.Cells.Item(41,3).Value2 = ' Best of luck! ♪'
$ws.Cells.Item(41,4).Value2 = ' Удачи вам всем! ♪'
$ws.Cells.Item(41,5).Value2 = ' Ôäàœé âàí âòåí! ♪'
$excel.CutCopyMode = $false

# Final selection to match target view state
$ws.Range("E41").Select()
